$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update personality / emotion text labels
$ws.Range("B2").Value = "Hate"
$ws.Range("D2").Value = "Fired"
$ws.Range("F3").Value = "High Conscientiousness"
$ws.Range("F7").Value = "Low Openness"

# Update numeric values
$ws.Range("A2").Value = -1.6607036590576172
$ws.Range("C2").Value = 5.347113132476807
